$wb = $excel.ActiveWorkbook

# --- 1. Selection / view tweaks on existing sheets ---

# token_recuperacao: was the active tab with a single-cell selection (D14);
# becomes non-active with the whole used range selected.
$wsToken = $wb.Worksheets.Item("token_recuperacao")
$wsToken.Range("A1:D7").Select()

# banco_de_leite: selection moves from the whole range to cell F38.
$wsBanco = $wb.Worksheets.Item("banco_de_leite")
$wsBanco.Activate()
$wsBanco.Range("F38").Select()

# Notificacao: selection moves from D30 to H23.
$wsNotif = $wb.Worksheets.Item("Notificacao")
$wsNotif.Activate()
$wsNotif.Range("H23").Select()

# --- 2. Rename Notificacao -> notificacao (lower case) ---
$wsNotif.Name = "notificacao"

# --- 3. Add the new "usuario_senha_historico" sheet at the end ---
# Copy token_recuperacao since it shares the exact same layout/styles
# (merged title row, header row, 4 data columns) that the new sheet needs.
$wsToken.Copy([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsNew = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsNew.Name = "usuario_senha_historico"

# --- 4. Populate the new sheet's content ---
$wsNew.Range("A1").Value = "Tabela: usuario_senha_historico"

$wsNew.Range("A2").Value = "Campo"
$wsNew.Range("B2").Value = "Tipo"
$wsNew.Range("C2").Value = "Tamanho"
$wsNew.Range("D2").Value = "Descrição"

$wsNew.Range("A3").Value = "id"
$wsNew.Range("B3").Value = "Serial"
$wsNew.Range("C3").Value = "-"
$wsNew.Range("D3").Value = "Chave primária da tabela bancos_de_leite"

$wsNew.Range("A4").Value = "usuario_id"
$wsNew.Range("B4").Value = "Integer"
$wsNew.Range("C4").Value = "-"
$wsNew.Range("D4").Value = "Usuario que esta solicitando a recuperação da senha"

$wsNew.Range("A5").Value = "senha"
$wsNew.Range("B5").Value = "Varchar"
$wsNew.Range("C5").Value = 6
$wsNew.Range("D5").Value = "senha de acesso ao sistema criptografada que foi alterada"

$wsNew.Range("A6").Value = "data_alteracao"
$wsNew.Range("B6").Value = "Timestamp"
$wsNew.Range("C6").Value = "-"
$wsNew.Range("D6").Value = "data de alteração da senha"

# the source sheet had a 7th data row that the new sheet doesn't need
$wsNew.Rows.Item(7).Delete()

# approximate the autofit width Excel would have computed for column D
$wsNew.Columns.Item(4).ColumnWidth = 50

# new sheet is the active one, with D5 selected
$wsNew.Range("D5").Select()

Write-Output "done"
